$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record (Femacal de La Calera - Achicoria) is inserted
# above the existing row 27, pushing every subsequent record down by one
# row (old row 122 becomes row 123, dimension grows from R122 to R123).
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row with the new record's data. Columns that
# are constant for every row in this subset (market/region/category/etc.)
# are filled in the same way as the surrounding rows.
$ws.Range("A27").Value = 3
$ws.Range("B27").Value = "Femacal de La Calera"
$ws.Range("C27").Value = "Coquimbo"
$ws.Range("D27").Value2 = 44459
$ws.Range("E27").Value = 5
$ws.Range("F27").Value = 100112010
$ws.Range("G27").Value = "Achicoria"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 105
$ws.Range("K27").Value = 5500
$ws.Range("L27").Value = 6000
$ws.Range("M27").Value = 5762
$ws.Range("N27").Value = "$/caja 16 unidades"
$ws.Range("O27").Value = "Provincia de Quillota"
$ws.Range("P27").Value = 360
$ws.Range("Q27").Value = 16
$ws.Range("R27").Value = "Hortaliza"
